$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1,1)
$c.Range.Text = "43-15="
$c = $t.Cell(1,2)
$c.Range.Text = "30-16="
$c = $t.Cell(1,3)
$c.Range.Text = "43-6="
$c = $t.Cell(1,4)
$c.Range.Text = "54+31="
$c = $t.Cell(1,5)
$c.Range.Text = "29+1="
$c = $t.Cell(2,1)
$c.Range.Text = "15-7="
$c = $t.Cell(2,2)
$c.Range.Text = "31+30="
$c = $t.Cell(2,3)
$c.Range.Text = "94-38="
$c = $t.Cell(2,4)
$c.Range.Text = "61-51="
$c = $t.Cell(2,5)
$c.Range.Text = "83-47="
$c = $t.Cell(3,1)
$c.Range.Text = "43+6="
$c = $t.Cell(3,2)
$c.Range.Text = "86-13="
$c = $t.Cell(3,3)
$c.Range.Text = "59+31="
$c = $t.Cell(3,4)
$c.Range.Text = "29-16="
$c = $t.Cell(3,5)
$c.Range.Text = "34-11="
$c = $t.Cell(4,1)
$c.Range.Text = "39+49="
$c = $t.Cell(4,2)
$c.Range.Text = "0+97="
$c = $t.Cell(4,3)
$c.Range.Text = "59+36="
$c = $t.Cell(4,4)
$c.Range.Text = "2+23="
$c = $t.Cell(4,5)
$c.Range.Text = "80-78="
$c = $t.Cell(5,1)
$c.Range.Text = "51-10="
$c = $t.Cell(5,2)
$c.Range.Text = "5+78="
$c = $t.Cell(5,3)
$c.Range.Text = "53-8="
$c = $t.Cell(5,4)
$c.Range.Text = "93-41="
$c = $t.Cell(5,5)
$c.Range.Text = "23+20="
$c = $t.Cell(6,1)
$c.Range.Text = "10+76="
$c = $t.Cell(6,2)
$c.Range.Text = "37+56="
$c = $t.Cell(6,3)
$c.Range.Text = "34+63="
$c = $t.Cell(6,4)
$c.Range.Text = "66-32="
$c = $t.Cell(6,5)
$c.Range.Text = "7+7="
$c = $t.Cell(7,1)
$c.Range.Text = "42-14="
$c = $t.Cell(7,2)
$c.Range.Text = "60-15="
$c = $t.Cell(7,3)
$c.Range.Text = "1+97="
$c = $t.Cell(7,4)
$c.Range.Text = "46+26="
$c = $t.Cell(7,5)
$c.Range.Text = "81-11="
$c = $t.Cell(8,1)
$c.Range.Text = "9+81="
$c = $t.Cell(8,2)
$c.Range.Text = "97-31="
$c = $t.Cell(8,3)
$c.Range.Text = "98-90="
$c = $t.Cell(8,4)
$c.Range.Text = "50+42="
$c = $t.Cell(8,5)
$c.Range.Text = "4+38="
$c = $t.Cell(9,1)
$c.Range.Text = "39+27="
$c = $t.Cell(9,2)
$c.Range.Text = "44+35="
$c = $t.Cell(9,3)
$c.Range.Text = "86-45="
$c = $t.Cell(9,4)
$c.Range.Text = "76-49="
$c = $t.Cell(9,5)
$c.Range.Text = "42+13="
$c = $t.Cell(10,1)
$c.Range.Text = "56+4="
$c = $t.Cell(10,2)
$c.Range.Text = "24+0="
$c = $t.Cell(10,3)
$c.Range.Text = "26-1="
$c = $t.Cell(10,4)
$c.Range.Text = "76+4="
$c = $t.Cell(10,5)
$c.Range.Text = "27-15="
$c = $t.Cell(11,1)
$c.Range.Text = "57-31="
$c = $t.Cell(11,2)
$c.Range.Text = "55+1="
$c = $t.Cell(11,3)
$c.Range.Text = "65-3="
$c = $t.Cell(11,4)
$c.Range.Text = "93-89="
$c = $t.Cell(11,5)
$c.Range.Text = "67+11="
$c = $t.Cell(12,1)
$c.Range.Text = "50-14="
$c = $t.Cell(12,2)
$c.Range.Text = "36-19="
$c = $t.Cell(12,3)
$c.Range.Text = "44-37="
$c = $t.Cell(12,4)
$c.Range.Text = "68-41="
$c = $t.Cell(12,5)
$c.Range.Text = "18+58="
$c = $t.Cell(13,1)
$c.Range.Text = "93-76="
$c = $t.Cell(13,2)
$c.Range.Text = "28+19="
$c = $t.Cell(13,3)
$c.Range.Text = "62+32="
$c = $t.Cell(13,4)
$c.Range.Text = "70-7="
$c = $t.Cell(13,5)
$c.Range.Text = "70-64="
$c = $t.Cell(14,1)
$c.Range.Text = "74-21="
$c = $t.Cell(14,2)
$c.Range.Text = "77-45="
$c = $t.Cell(14,3)
$c.Range.Text = "92-50="
$c = $t.Cell(14,4)
$c.Range.Text = "3+96="
$c = $t.Cell(14,5)
$c.Range.Text = "30+32="
$c = $t.Cell(15,1)
$c.Range.Text = "88-59="
$c = $t.Cell(15,2)
$c.Range.Text = "89-72="
$c = $t.Cell(15,3)
$c.Range.Text = "9-2="
$c = $t.Cell(15,4)
$c.Range.Text = "1+85="
$c = $t.Cell(15,5)
$c.Range.Text = "77-4="
$c = $t.Cell(16,1)
$c.Range.Text = "67+27="
$c = $t.Cell(16,2)
$c.Range.Text = "94-22="
$c = $t.Cell(16,3)
$c.Range.Text = "74-11="
$c = $t.Cell(16,4)
$c.Range.Text = "31+51="
$c = $t.Cell(16,5)
$c.Range.Text = "84-36="
$c = $t.Cell(17,1)
$c.Range.Text = "73-23="
$c = $t.Cell(17,2)
$c.Range.Text = "44+50="
$c = $t.Cell(17,3)
$c.Range.Text = "75+5="
$c = $t.Cell(17,4)
$c.Range.Text = "92-62="
$c = $t.Cell(17,5)
$c.Range.Text = "55+26="
$c = $t.Cell(18,1)
$c.Range.Text = "93-88="
$c = $t.Cell(18,2)
$c.Range.Text = "92-78="
$c = $t.Cell(18,3)
$c.Range.Text = "12+21="
$c = $t.Cell(18,4)
$c.Range.Text = "46+31="
$c = $t.Cell(18,5)
$c.Range.Text = "22+77="
$c = $t.Cell(19,1)
$c.Range.Text = "35-4="
$c = $t.Cell(19,2)
$c.Range.Text = "90-19="
$c = $t.Cell(19,3)
$c.Range.Text = "72-1="
$c = $t.Cell(19,4)
$c.Range.Text = "24+61="
$c = $t.Cell(19,5)
$c.Range.Text = "88+1="
$c = $t.Cell(20,1)
$c.Range.Text = "79-76="
$c = $t.Cell(20,2)
$c.Range.Text = "75-13="
$c = $t.Cell(20,3)
$c.Range.Text = "67-38="
$c = $t.Cell(20,4)
$c.Range.Text = "57+12="
$c = $t.Cell(20,5)
$c.Range.Text = "0+86="
